{"js": "// Apply the diff: update the date line and the 25 division problems.\n// Each \"find\" string below is unique in the document, so searching for\n// the original text and replacing it in place is safe even though some\n// replacement values happen to equal other cells' original values\n// (e.g. \"88\u00f77=\" -> \"37\u00f78=\" while a different cell's original \"37\u00f78=\"\n// -> \"50\u00f73=\"): we search for the literal old text for every pair\n// against the (still-original) document content, so there is no\n// cross-talk between pairs.\nconst replacements = [\n  [\"2023-12-31 Sunday\", \"2024-01-01 Monday\"],\n  [\"71\u00f74=\", \"69\u00f77=\"],\n  [\"19\u00f78=\", \"39\u00f74=\"],\n  [\"36\u00f76=\", \"13\u00f74=\"],\n  [\"29\u00f76=\", \"25\u00f73=\"],\n  [\"80\u00f78=\", \"77\u00f79=\"],\n  [\"63\u00f79=\", \"54\u00f72=\"],\n  [\"68\u00f75=\", \"93\u00f72=\"],\n  [\"20\u00f78=\", \"44\u00f74=\"],\n  [\"37\u00f78=\", \"50\u00f73=\"],\n  [\"99\u00f79=\", \"77\u00f74=\"],\n  [\"69\u00f73=\", \"16\u00f79=\"],\n  [\"42\u00f74=\", \"33\u00f72=\"],\n  [\"30\u00f74=\", \"98\u00f75=\"],\n  [\"87\u00f78=\", \"62\u00f73=\"],\n  [\"89\u00f75=\", \"62\u00f73=\"],\n  [\"91\u00f77=\", \"77\u00f76=\"],\n  [\"88\u00f77=\", \"37\u00f78=\"],\n  [\"61\u00f72=\", \"51\u00f73=\"],\n  [\"34\u00f72=\", \"83\u00f73=\"],\n  [\"79\u00f75=\", \"44\u00f79=\"],\n  [\"86\u00f76=\", \"79\u00f79=\"],\n  [\"86\u00f73=\", \"96\u00f76=\"],\n  [\"59\u00f73=\", \"20\u00f72=\"],\n  [\"14\u00f76=\", \"35\u00f78=\"],\n  [\"92\u00f76=\", \"83\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the diff: update the date line and the 25 division problems.\n# Every \"old\" string below is unique in the document, so a targeted\n# Find/Replace (MatchCase, whole text, not wildcard, Replace:=1 i.e.\n# wdReplaceOne) against the full document content is safe for every\n# pair -- including \"88\u00f77=\" -> \"37\u00f78=\" together with the unrelated\n# \"37\u00f78=\" -> \"50\u00f73=\" pair -- because each Find.Execute only ever\n# matches the literal original text of its own pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ old = \"2023-12-31 Sunday\"; new = \"2024-01-01 Monday\" },\n    @{ old = \"71\u00f74=\"; new = \"69\u00f77=\" },\n    @{ old = \"19\u00f78=\"; new = \"39\u00f74=\" },\n    @{ old = \"36\u00f76=\"; new = \"13\u00f74=\" },\n    @{ old = \"29\u00f76=\"; new = \"25\u00f73=\" },\n    @{ old = \"80\u00f78=\"; new = \"77\u00f79=\" },\n    @{ old = \"63\u00f79=\"; new = \"54\u00f72=\" },\n    @{ old = \"68\u00f75=\"; new = \"93\u00f72=\" },\n    @{ old = \"20\u00f78=\"; new = \"44\u00f74=\" },\n    @{ old = \"37\u00f78=\"; new = \"50\u00f73=\" },\n    @{ old = \"99\u00f79=\"; new = \"77\u00f74=\" },\n    @{ old = \"69\u00f73=\"; new = \"16\u00f79=\" },\n    @{ old = \"42\u00f74=\"; new = \"33\u00f72=\" },\n    @{ old = \"30\u00f74=\"; new = \"98\u00f75=\" },\n    @{ old = \"87\u00f78=\"; new = \"62\u00f73=\" },\n    @{ old = \"89\u00f75=\"; new = \"62\u00f73=\" },\n    @{ old = \"91\u00f77=\"; new = \"77\u00f76=\" },\n    @{ old = \"88\u00f77=\"; new = \"37\u00f78=\" },\n    @{ old = \"61\u00f72=\"; new = \"51\u00f73=\" },\n    @{ old = \"34\u00f72=\"; new = \"83\u00f73=\" },\n    @{ old = \"79\u00f75=\"; new = \"44\u00f79=\" },\n    @{ old = \"86\u00f76=\"; new = \"79\u00f79=\" },\n    @{ old = \"86\u00f73=\"; new = \"96\u00f76=\" },\n    @{ old = \"59\u00f73=\"; new = \"20\u00f72=\" },\n    @{ old = \"14\u00f76=\"; new = \"35\u00f78=\" },\n    @{ old = \"92\u00f76=\"; new = \"83\u00f74=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 1)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $($pair.old)\"\n    }\n}\n\n$d.Save()\n"}
